$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column widths (Excel pads ColumnWidth by ~0.8333 chars when it stores the
# saved <col width="..."/>, so we dial the assigned width back by that amount to
# land exactly on the target stored widths of 15/15/25/15/15/10/10/10).
$ws.Columns.Item(1).ColumnWidth = 14.166666666666666
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(3).ColumnWidth = 24.166666666666668
$ws.Columns.Item(4).ColumnWidth = 14.166666666666666
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666
$ws.Columns.Item(6).ColumnWidth = 9.166666666666666
$ws.Columns.Item(7).ColumnWidth = 9.166666666666666
$ws.Columns.Item(8).ColumnWidth = 9.166666666666666

# Update header row, translated to Polish and extended with new columns
$ws.Range("A1").Value = "UŻYTKOWNIK"
$ws.Range("B1").Value = "DATA WNIOSKU"
$ws.Range("C1").Value = "TYP"
$ws.Range("D1").Value = "OD"
$ws.Range("E1").Value = "DO"
$ws.Range("F1").Value = "ILOŚĆ DNI"
$ws.Range("G1").Value = "AKCEPTACJA"

# Make sure the date-like text in E:F keeps being stored as plain text
# (matches the source workbook, where dates are strings, not real dates)
# -- without this, Excel auto-converts "2023-06-20"-style text into a date serial.
$ws.Range("E2:F5").NumberFormat = "@"

# Row 2 - Alex Admin: 2023-06-20 -> 2023-07-09, 17 days, accepted
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = "wypoczynkowy"
$ws.Range("E2").Value = "2023-06-20"
$ws.Range("F2").Value = "2023-07-09"
$ws.Range("G2").Value = 17
$ws.Range("H2").Value = "Tak"

# Row 3 - Noll Roman: 2023-06-20 -> 2023-06-21, 2 days, accepted
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = "wypoczynkowy"
$ws.Range("E3").Value = "2023-06-20"
$ws.Range("F3").Value = "2023-06-21"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = "Tak"

# Row 4 - Noll Roman: 2023-06-20 -> 2023-06-22, 3 days, denied
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = "wypoczynkowy"
$ws.Range("E4").Value = "2023-06-20"
$ws.Range("F4").Value = "2023-06-22"
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = "Nie"

# Row 5 - Noll Roman: 2023-06-21 -> 2023-06-22, 2 days, denied
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = "wypoczynkowy"
$ws.Range("E5").Value = "2023-06-21"
$ws.Range("F5").Value = "2023-06-22"
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = "Nie"
